{"js": "// Replace each three-digit-by-one-digit multiplication answer with its\n// updated version. Each \"old\" string is a unique, literal run of text\n// (e.g. \"402\u00d72=804\") living inside a single <w:t>, so an exact,\n// case-sensitive search-and-replace on context.document.body locates the\n// correct cell even though a couple of the new values momentarily equal an\n// old value used elsewhere in the table (the pairs are applied strictly in\n// document order, so that never collides).\nconst replacements = [\n  { old: \"402\u00d72=804\", new: \"543\u00d72=1086\" },\n  { old: \"129\u00d74=516\", new: \"102\u00d72=204\" },\n  { old: \"970\u00d77=6790\", new: \"628\u00d75=3140\" },\n  { old: \"874\u00d77=6118\", new: \"322\u00d79=2898\" },\n  { old: \"773\u00d77=5411\", new: \"408\u00d76=2448\" },\n  { old: \"168\u00d72=336\", new: \"692\u00d75=3460\" },\n  { old: \"443\u00d79=3987\", new: \"913\u00d77=6391\" },\n  { old: \"828\u00d74=3312\", new: \"788\u00d77=5516\" },\n  { old: \"757\u00d75=3785\", new: \"831\u00d76=4986\" },\n  { old: \"943\u00d73=2829\", new: \"342\u00d76=2052\" },\n  { old: \"947\u00d74=3788\", new: \"304\u00d75=1520\" },\n  { old: \"465\u00d75=2325\", new: \"963\u00d77=6741\" },\n  { old: \"959\u00d79=8631\", new: \"125\u00d77=875\" },\n  { old: \"693\u00d73=2079\", new: \"269\u00d75=1345\" },\n  { old: \"772\u00d76=4632\", new: \"947\u00d74=3788\" },\n  { old: \"796\u00d78=6368\", new: \"197\u00d78=1576\" },\n  { old: \"324\u00d73=972\", new: \"189\u00d78=1512\" },\n  { old: \"902\u00d76=5412\", new: \"609\u00d77=4263\" },\n  { old: \"871\u00d75=4355\", new: \"774\u00d79=6966\" },\n  { old: \"188\u00d78=1504\", new: \"942\u00d75=4710\" },\n  { old: \"551\u00d75=2755\", new: \"404\u00d76=2424\" },\n  { old: \"560\u00d76=3360\", new: \"398\u00d79=3582\" },\n  { old: \"376\u00d74=1504\", new: \"697\u00d76=4182\" },\n  { old: \"545\u00d75=2725\", new: \"173\u00d72=346\" },\n  { old: \"354\u00d77=2478\", new: \"869\u00d75=4345\" },\n];\n\nconst body = context.document.body;\n\nfor (const { old, new: replacement } of replacements) {\n  const results = body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${old}`);\n  }\n\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication answer with its\n# updated version. Each \"old\" string is a unique, literal run of text\n# (e.g. \"402x2=804\") found inside exactly one table cell, so a single\n# Find/Replace (wdReplaceOne) against the whole document body locates the\n# right cell. The pairs are applied strictly in the order they occur in\n# the document -- this matters because a couple of the \"new\" values\n# momentarily equal an \"old\" value used further down the table, and doing\n# the replacements top-to-bottom (re-searching $d.Content, i.e. from the\n# very start, each time) guarantees every Find lands on the original,\n# not-yet-updated cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"402\u00d72=804\";   New = \"543\u00d72=1086\" },\n    @{ Old = \"129\u00d74=516\";   New = \"102\u00d72=204\" },\n    @{ Old = \"970\u00d77=6790\";  New = \"628\u00d75=3140\" },\n    @{ Old = \"874\u00d77=6118\";  New = \"322\u00d79=2898\" },\n    @{ Old = \"773\u00d77=5411\";  New = \"408\u00d76=2448\" },\n    @{ Old = \"168\u00d72=336\";   New = \"692\u00d75=3460\" },\n    @{ Old = \"443\u00d79=3987\";  New = \"913\u00d77=6391\" },\n    @{ Old = \"828\u00d74=3312\";  New = \"788\u00d77=5516\" },\n    @{ Old = \"757\u00d75=3785\";  New = \"831\u00d76=4986\" },\n    @{ Old = \"943\u00d73=2829\";  New = \"342\u00d76=2052\" },\n    @{ Old = \"947\u00d74=3788\";  New = \"304\u00d75=1520\" },\n    @{ Old = \"465\u00d75=2325\";  New = \"963\u00d77=6741\" },\n    @{ Old = \"959\u00d79=8631\";  New = \"125\u00d77=875\" },\n    @{ Old = \"693\u00d73=2079\";  New = \"269\u00d75=1345\" },\n    @{ Old = \"772\u00d76=4632\";  New = \"947\u00d74=3788\" },\n    @{ Old = \"796\u00d78=6368\";  New = \"197\u00d78=1576\" },\n    @{ Old = \"324\u00d73=972\";   New = \"189\u00d78=1512\" },\n    @{ Old = \"902\u00d76=5412\";  New = \"609\u00d77=4263\" },\n    @{ Old = \"871\u00d75=4355\";  New = \"774\u00d79=6966\" },\n    @{ Old = \"188\u00d78=1504\";  New = \"942\u00d75=4710\" },\n    @{ Old = \"551\u00d75=2755\";  New = \"404\u00d76=2424\" },\n    @{ Old = \"560\u00d76=3360\";  New = \"398\u00d79=3582\" },\n    @{ Old = \"376\u00d74=1504\";  New = \"697\u00d76=4182\" },\n    @{ Old = \"545\u00d75=2725\";  New = \"173\u00d72=346\" },\n    @{ Old = \"354\u00d77=2478\";  New = \"869\u00d75=4345\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    # 0 = wdFindStop, 1 = wdReplaceOne -- replace only the first (and only)\n    # match, not every occurrence.\n    $found = $find.Execute($null, $false, $false, $false, $false, $false, $true, 0, $false, $r.New, 1)\n    if (-not $found) {\n        throw \"Could not find text to replace: $($r.Old)\"\n    }\n}\n"}
